# Initial check-in of translations changes.
#
# The "display.*" vocabulary used throughout the farm_crop survey/settings
# sheets is renamed to its ".text"-suffixed form, and the now-unused
# "display.new_instance_text" column on the survey sheet is dropped.
#   display.text               -> display.prompt.text
#   display.hint                -> display.hint.text
#   display.title               -> display.title.text
#   display.new_instance_text   -> (removed column)

$wb = $excel.ActiveWorkbook

# --- settings sheet: display.title -> display.title.text -------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("C1").Value = "display.title.text"

# --- initial sheet: display.text -> display.prompt.text --------------------
$initial = $wb.Worksheets.Item("initial")
$initial.Range("C1").Value = "display.prompt.text"

# --- survey sheet: display.text/display.hint -> *.text, drop new_instance --
$survey = $wb.Worksheets.Item("survey")
$survey.Range("F1").Value = "display.prompt.text"
$survey.Range("G1").Value = "display.hint.text"
$survey.Columns.Item(9).Delete()

# --- restore per-sheet selections -------------------------------------------
$initial.Range("E3").Select() | Out-Null
$survey.Range("C9").Select() | Out-Null
$settings.Range("C2").Select() | Out-Null

# --- the active tab moves from "model" to "settings" ------------------------
$settings.Activate()
